$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures as literal text (the sheet is generated by a scraper,
# not typed in by a user), including values that look numeric, e.g. "9.70", "0.999",
# "0.0000250". A plain `.Value = "9.70"` assignment would make Excel infer a Number
# and silently drop the significant trailing/leading zeros. Prefixing the literal with
# a single quote forces Excel to keep/store it as Text (quoted-text cell entry),
# exactly reproducing the original inlineStr text cells.

# Row 2
$ws.Range("D2").Value = "'69.031.46"
$ws.Range("E2").Value = "  +2.30%  "

# Row 3
$ws.Range("D3").Value = "'3.818.72"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "'629.88"
$ws.Range("E5").Value = "  +5.45%  "

# Row 6
$ws.Range("D6").Value = "'165.24"
$ws.Range("E6").Value = "  +0.55%  "

# Row 7
$ws.Range("D7").Value = "'3.816.53"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("E9").Value = "  +1.13%  "

# Row 10
$ws.Range("E10").Value = "  +2.52%  "

# Row 11
$ws.Range("E11").Value = "  +1.26%  "

# Row 12
$ws.Range("D12").Value = "'6.61"
$ws.Range("E12").Value = "  +3.34%  "

# Row 13
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +1.17%  "

# Row 14
$ws.Range("D14").Value = "'36.03"
$ws.Range("E14").Value = "  +1.35%  "

# Row 15
$ws.Range("D15").Value = "'4.458.57"
$ws.Range("E15").Value = "  +0.96%  "

# Row 16
$ws.Range("D16").Value = "'3.745.10"
$ws.Range("E16").Value = "  -0.35%  "

# Row 17
$ws.Range("D17").Value = "'69.067.09"
$ws.Range("E17").Value = "  +2.21%  "

# Row 18
$ws.Range("D18").Value = "'17.99"
$ws.Range("E18").Value = "  -1.30%  "

# Row 19
$ws.Range("E19").Value = "  +1.87%  "

# Row 20
$ws.Range("E20").Value = "  +0.08%  "

# Row 21
$ws.Range("D21").Value = "'465.17"
$ws.Range("E21").Value = "  +1.10%  "

# Row 22
$ws.Range("D22").Value = "'9.70"

# Row 23
$ws.Range("E23").Value = "  +1.93%  "

# Row 24
$ws.Range("E24").Value = "  +4.71%  "

# Row 25
$ws.Range("D25").Value = "'83.72"
$ws.Range("E25").Value = "  +1.71%  "

# Row 26
$ws.Range("E26").Value = "  +0.10%  "

# Row 27
$ws.Range("E27").Value = "  +3.39%  "

# Row 28
$ws.Range("E28").Value = "  +0.68%  "

# Row 29
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("D30").Value = "'3.966.76"
$ws.Range("E30").Value = "  +0.99%  "

# Row 31
$ws.Range("D31").Value = "'2.70"
$ws.Range("E31").Value = "  +1.97%  "

# Row 32
$ws.Range("E32").Value = "  +1.92%  "

# Row 33
$ws.Range("D33").Value = "'7.28"
$ws.Range("E33").Value = "  -1.50%  "

# Row 34
$ws.Range("D34").Value = "'29.22"
$ws.Range("E34").Value = "  +1.15%  "

# Row 35
$ws.Range("E35").Value = "  +0.40%  "

# Row 36
$ws.Range("E36").Value = "  +1.51%  "

# Row 37
$ws.Range("E37").Value = "  +3.23%  "

# Row 38
$ws.Range("E38").Value = "  +7.98%  "

# Row 39
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  +5.95%  "

# Row 40
$ws.Range("E40").Value = "  +3.35%  "

# Row 41
$ws.Range("E41").Value = "  -0.80%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 44
$ws.Range("D44").Value = "'157.27"
$ws.Range("E44").Value = "  +4.22%  "

# Row 45
$ws.Range("E45").Value = "  +5.64%  "

# Row 46
$ws.Range("E46").Value = "  +1.30%  "

# Row 47
$ws.Range("D47").Value = "'46.78"
$ws.Range("E47").Value = "  -1.45%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'8.45"
$ws.Range("E48").Value = "  +1.85%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.90"
$ws.Range("E49").Value = "  +3.24%  "

# Row 50
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'42.43"
$ws.Range("E50").Value = "  -2.65%  "

# Row 51
$ws.Range("D51").Value = "'0.000279"
$ws.Range("E51").Value = "  +13.59%  "
